$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing row formats before we shuffle values around, so the
# "border10" (mid-list) and "border11" (tail-of-list) looks can be re-applied
# to the right rows once the new, longer scenario list is written out.
$midFormatRange = $ws.Range("A2")
$tailFormatRange = $ws.Range("A5")

# New scenario list (replaces the old 4-item list with the refreshed one:
# "E2E_BTS_Bundle" is gone, "E2E_BTS_Service" + four more scenarios added).
$ws.Range("A1").Value = "Scenarios"
$ws.Range("A2").Value = "E2E_24_PS_ES_RE_Data"
$ws.Range("A3").Value = "E2E_Intra"
$ws.Range("A4").Value = "E2E_BTS_Service"
$ws.Range("A5").Value = "E2E_CTO_BTS"
$ws.Range("A6").Value = "E2E_23_RENEWAL_1"
$ws.Range("A7").Value = "E2E_23_RENEWAL_2"
$ws.Range("A8").Value = "E2E_Installments"
$ws.Range("A9").Value = "E2E_20_RFC2"
$ws.Range("A10").Value = "E2E_StockRotationReturnDelivery"

# Re-apply the two row styles to the (now larger) ranges they cover.
$midFormatRange.Copy() | Out-Null
$ws.Range("A2:A4").PasteSpecial(-4122) | Out-Null

$tailFormatRange.Copy() | Out-Null
$ws.Range("A5:A10").PasteSpecial(-4122) | Out-Null

$null = $ws.Range("G15").Select()
